# Remove the notion of pre-installed applications in the overview slide.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1. Delete the "Application 1" rounded rectangle.
$s.Shapes.Item("Rounded Rectangle 10").Delete()

# 2. Merge the "EE " + "Port" runs of the "VEE Port" shape into a single
#    run "EE Port" (text itself is unchanged - only run structure).
$sh = $s.Shapes.Item("Rounded Rectangle 24")
$tr = $sh.TextFrame.TextRange
$sub = $tr.Characters(2, 7)
$sub.Text = "EE Port"

# 3. Delete the "Application N" rounded rectangle and the down arrow that
#    pointed to it ("Down Arrow 30"); the remaining "Down Arrow 31" shape
#    slides into that slot.
$s.Shapes.Item("Rounded Rectangle 71").Delete()
$s.Shapes.Item("Down Arrow 30").Delete()

# 4. Delete the dashed bounding rectangle and its
#    "Pre-installed Applications" caption.
$s.Shapes.Item("Rectangle 21").Delete()
$s.Shapes.Item("Rectangle 22").Delete()

# 5. Merge the "Kernel " + "binary" runs of the kernel-binary shape into a
#    single run "Kernel binary" (text itself is unchanged).
$sh2 = $s.Shapes.Item("Rounded Rectangle 23")
$tr2 = $sh2.TextFrame.TextRange
$sub2 = $tr2.Characters(1, 13)
$sub2.Text = "Kernel binary"
